# Aula 07 Microcontroladores ajuste 17maio2023
# Slide 2, title placeholder: split the trailing run
#   "Periféricos Externos – Sensores e Atuadores"
# into
#   "Periféricos – " + "Sensores e Atuadores"
# (i.e. remove "Externos " right after "Periféricos ").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

# Locate the run's text within the whole paragraph text so the edit is
# anchored on content rather than hard-coded offsets.
$full   = $tr.Text
$oldSeg = "Periféricos Externos – Sensores e Atuadores"
$startIdx = $full.IndexOf($oldSeg)
if ($startIdx -lt 0) {
    throw "Could not find target text in title placeholder"
}

$headOld = "Periféricos Externos – "
$headNew = "Periféricos – "

# 1-based character index where the run of interest starts.
$start = $startIdx + 1

# Replace just the "Periféricos Externos – " portion with "Periféricos – ",
# which splits the original single run into two runs: the edited head and
# the untouched tail ("Sensores e Atuadores").
$head = $tr.Characters($start, $headOld.Length)
$head.Text = $headNew
